$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap store names between row 4 and row 5 (Manauara <-> Ponta Negra)
$ws.Range("A4").Value = "Bibi Cell Ponta Negra"
$ws.Range("A5").Value = "Bibi Cell Manauara"

# Row 2: update AB2, AC2, AG2
$ws.Range("AB2").Value = 10793.22
$ws.Range("AC2").Value = 1817
$ws.Range("AG2").Value = 284807.17

# Row 3: update AB3, AC3, AG3
$ws.Range("AB3").Value = 3858.4
$ws.Range("AC3").Value = 4344
$ws.Range("AG3").Value = 149960.99

# Row 4: swap data with row 5 (cols B..AA) and set new AB4, AC4, AD4, AG4
$ws.Range("B4").Value = 1800.01
$ws.Range("C4").Value = 4670
$ws.Range("D4").Value = 1748.51
$ws.Range("E4").Value = 5592
$ws.Range("F4").Value = 3002
$ws.Range("G4").Value = 823
$ws.Range("H4").Value = 3138.5
$ws.Range("I4").Value = 1613
$ws.Range("J4").Value = 2786.02
$ws.Range("K4").Value = 6097.5
$ws.Range("L4").Value = 3514.36
$ws.Range("M4").Value = 5434.87
$ws.Range("N4").Value = 1478
$ws.Range("O4").Value = 4390.5
$ws.Range("P4").Value = 1481.42
$ws.Range("Q4").Value = 2772
$ws.Range("R4").Value = 1781
$ws.Range("S4").Value = 1114
$ws.Range("T4").Value = 1142.5
$ws.Range("U4").Value = 1795
$ws.Range("V4").Value = 2338.01
$ws.Range("W4").Value = 1118.5
$ws.Range("X4").Value = 3358.8
$ws.Range("Y4").Value = 2671
$ws.Range("Z4").Value = 2711.11
$ws.Range("AA4").Value = 6123.4
$ws.Range("AB4").Value = 7033.9
$ws.Range("AC4").Value = 2211
$ws.Range("AD4").Value = 4540.8
$ws.Range("AG4").Value = 88280.71000000001

# Row 5: swap data with row 4 (cols B..AA) and set new AB5, AC5, AD5, AG5
$ws.Range("B5").Value = 3340
$ws.Range("C5").Value = 1519
$ws.Range("D5").Value = 2934
$ws.Range("E5").Value = 1819
$ws.Range("F5").Value = 2503
$ws.Range("G5").Value = 2892
$ws.Range("H5").Value = 4208.4
$ws.Range("I5").Value = 3329.9
$ws.Range("J5").Value = 4038
$ws.Range("K5").Value = 2830.9
$ws.Range("L5").Value = 3525
$ws.Range("M5").Value = 2422
$ws.Range("N5").Value = 2493
$ws.Range("O5").Value = 5411
$ws.Range("P5").Value = 3140
$ws.Range("Q5").Value = 2599
$ws.Range("R5").Value = 2746
$ws.Range("S5").Value = 3199.9
$ws.Range("T5").Value = 3780
$ws.Range("U5").Value = 3561.5
$ws.Range("V5").Value = 4526
$ws.Range("W5").Value = 2936
$ws.Range("X5").Value = 2117
$ws.Range("Y5").Value = 4144
$ws.Range("Z5").Value = 1379.75
$ws.Range("AA5").Value = 2569
$ws.Range("AB5").Value = 2834
$ws.Range("AC5").Value = 2531.9
$ws.Range("AD5").Value = 1589
$ws.Range("AG5").Value = 86918.25

# Row 6: update AB6, AC6, AD6, AG6
$ws.Range("AB6").Value = 24519.52
$ws.Range("AC6").Value = 10903.9
$ws.Range("AD6").Value = 6129.8
$ws.Range("AG6").Value = 609967.12
